$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("g18.3")

$ws.Range("B22").Value = 22.63548882549046
$ws.Range("C22").Value = 11.48481906767816
$ws.Range("D22").Value = 13.51093904698167

$ws.Range("B23").Value = 23.09893510493992
$ws.Range("C23").Value = 11.83471890624934
$ws.Range("D23").Value = 13.74061011986208

$ws.Range("B24").Value = 23.52577827880783
$ws.Range("C24").Value = 12.16101662962668
$ws.Range("D24").Value = 13.98682692633936

$ws.Range("B25").Value = 23.71146373109548
$ws.Range("C25").Value = 12.4483240726964
$ws.Range("D25").Value = 14.01600793116783

$ws.Range("B26").Value = 23.90207758604287
$ws.Range("C26").Value = 12.61414502123088
$ws.Range("D26").Value = 14.19608567023208
